# Update "想去人数" (want-to-go count) values in column F across the four
# sheets of the workbook, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 37655
$ws1.Range("F5").Value  = 767
$ws1.Range("F6").Value  = 481
$ws1.Range("F7").Value  = 366
$ws1.Range("F8").Value  = 466
$ws1.Range("F9").Value  = 843
$ws1.Range("F10").Value = 96
$ws1.Range("F11").Value = 716
$ws1.Range("F12").Value = 548
$ws1.Range("F13").Value = 46
$ws1.Range("F17").Value = 176
$ws1.Range("F20").Value = 1169
$ws1.Range("F22").Value = 830
$ws1.Range("F23").Value = 2532
$ws1.Range("F24").Value = 1010
$ws1.Range("F27").Value = 1163
$ws1.Range("F29").Value = 777
$ws1.Range("F30").Value = 59
$ws1.Range("F31").Value = 1160

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 402
$ws2.Range("F12").Value = 10

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 631

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 631
$ws4.Range("F3").Value  = 37655
$ws4.Range("F6").Value  = 767
$ws4.Range("F7").Value  = 481
$ws4.Range("F9").Value  = 366
$ws4.Range("F10").Value = 466
$ws4.Range("F11").Value = 402
$ws4.Range("F15").Value = 843
$ws4.Range("F16").Value = 96
$ws4.Range("F17").Value = 716
$ws4.Range("F18").Value = 548
$ws4.Range("F19").Value = 46
$ws4.Range("F25").Value = 19
$ws4.Range("F28").Value = 176
$ws4.Range("F31").Value = 1169
$ws4.Range("F33").Value = 830
$ws4.Range("F34").Value = 2532
$ws4.Range("F35").Value = 1010
$ws4.Range("F38").Value = 1163
$ws4.Range("F40").Value = 10
$ws4.Range("F41").Value = 777
$ws4.Range("F42").Value = 59
$ws4.Range("F43").Value = 1160

$wb.Save()
